# "added 4wk low sales check"
# Update the forecast figures (MyForecast, Inventory Coverage, Stockout Risk,
# Reorder Urgency, Seasonality Index) on the "Forecast Comparison" sheet and
# the corresponding roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$forecast = $wb.Worksheets.Item("Forecast Comparison")
$summary  = $wb.Worksheets.Item("Summary")

# --- MyForecast (column D) ---------------------------------------------
$forecast.Range("D2").Value  = 19
$forecast.Range("D3").Value  = 18
$forecast.Range("D4").Value  = 15
$forecast.Range("D5").Value  = 14
$forecast.Range("D6").Value  = 13
$forecast.Range("D7").Value  = 13
$forecast.Range("D8").Value  = 14
$forecast.Range("D9").Value  = 13
$forecast.Range("D11").Value = 9
$forecast.Range("D12").Value = 7
$forecast.Range("D13").Value = 6
$forecast.Range("D14").Value = 6
$forecast.Range("D15").Value = 6
$forecast.Range("D16").Value = 6
$forecast.Range("D17").Value = 6

# --- Inventory Coverage (column H) --------------------------------------
$forecast.Range("H2").Value  = 9.529999999999999
$forecast.Range("H3").Value  = 9
$forecast.Range("H4").Value  = 9.6
$forecast.Range("H5").Value  = 9.210000000000001
$forecast.Range("H6").Value  = 8.85
$forecast.Range("H7").Value  = 7.85
$forecast.Range("H8").Value  = 6.36
$forecast.Range("H9").Value  = 5.77
$forecast.Range("H10").Value = 5.17
$forecast.Range("H11").Value = 5.56
$forecast.Range("H12").Value = 5.86
$forecast.Range("H13").Value = 5.67
$forecast.Range("H14").Value = 4.67
$forecast.Range("H15").Value = 3.67
$forecast.Range("H16").Value = 2.67
$forecast.Range("H17").Value = 1.67

# --- Stockout Risk (column I) -------------------------------------------
$forecast.Range("I16").Value = "Low"
$forecast.Range("I17").Value = "Low"

# --- Reorder Urgency (column J) -----------------------------------------
$forecast.Range("J15").Value = "Normal"
$forecast.Range("J16").Value = "Normal"
$forecast.Range("J17").Value = "Normal"

# --- Seasonality Index (column L) ---------------------------------------
$forecast.Range("L2").Value  = 1.01
$forecast.Range("L3").Value  = 0.95
$forecast.Range("L4").Value  = 0.96
$forecast.Range("L5").Value  = 1.09
$forecast.Range("L6").Value  = 1.1
$forecast.Range("L7").Value  = 1.01
$forecast.Range("L8").Value  = 1.2
$forecast.Range("L9").Value  = 1.05
$forecast.Range("L10").Value = 0.82
$forecast.Range("L11").Value = 0.86
$forecast.Range("L13").Value = 1.06
$forecast.Range("L14").Value = 1.01
$forecast.Range("L15").Value = 0.83
$forecast.Range("L16").Value = 0.93

# --- Summary sheet roll-up totals ---------------------------------------
# Leading apostrophe keeps these as text (matching the rest of the
# "Value" column on this sheet) instead of Excel auto-converting the
# numeric-looking string to a number.
$summary.Range("B9").Value  = "'177"
$summary.Range("B10").Value = "'119"
$summary.Range("B11").Value = "'66"
$summary.Range("B12").Value = "'19"
$summary.Range("B14").Value = "'6"
